$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scenario rows appended to the results sheet (PositionsSteps_Tuba
# positive and negative tests, plus a login scenario), rows 121-149.
$rows = @(
    @("Login with valid username and password", "FAILED", "chrome"),
    @("Login with valid username and password", "FAILED", "chrome"),
    @("Add New Bank Accounts", "FAILED", "chrome"),
    @("Edit The Bank Accounts", "FAILED", "chrome"),
    @("Delete The Bank Accounts", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "UNDEFINED", "chrome"),
    @("Adding New Positions to the Admin Panel", "UNDEFINED", "chrome"),
    @("Adding New Positions to the Admin Panel", "UNDEFINED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "PASSED", "chrome"),
    @("Adding New Positions to the Admin Panel", "PASSED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "FAILED", "chrome"),
    @("Adding New Positions to the Admin Panel", "PASSED", "chrome")
)

$startRow = 121
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}
